# Re-number the "Раздел" (section) column (A) on the question sheet so
# each quiz section has 3 questions instead of 10 (rows 2-6 / section 1
# are left untouched, matching the source diff which only touches A7:A61).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sections = @(
    @{Start=7;  End=11; Value=2},
    @{Start=12; End=14; Value=3},
    @{Start=15; End=17; Value=4},
    @{Start=18; End=20; Value=5},
    @{Start=21; End=23; Value=6},
    @{Start=24; End=26; Value=7},
    @{Start=27; End=29; Value=8},
    @{Start=30; End=32; Value=9},
    @{Start=33; End=35; Value=10},
    @{Start=36; End=38; Value=11},
    @{Start=39; End=41; Value=12},
    @{Start=42; End=44; Value=13},
    @{Start=45; End=47; Value=14},
    @{Start=48; End=50; Value=15},
    @{Start=51; End=53; Value=16},
    @{Start=54; End=56; Value=17},
    @{Start=57; End=59; Value=18},
    @{Start=60; End=61; Value=19}
)

foreach ($section in $sections) {
    $range = $ws.Range("A$($section.Start):A$($section.End)")
    $range.Value = $section.Value
}

# Update the saved cursor/selection position to match the edited workbook
# (scrolled a bit further up, active cell moved to B48).
$ws.Range("B48").Select()

Write-Output "Section numbering updated; selection set to B48."
